# Revert "Powerpoint writer: consolidate text run nodes."
#
# The caption textbox on slide 1 currently stores its text as three
# runs: "The ", "picture ", "first". The target splits the trailing
# space off of the first two runs into their own runs, producing:
# "The", " ", "picture", " ", "first" -- five runs in total, with no
# change to the visible text.
#
# Re-assigning the Text of a Characters() sub-range that covers only
# part of an existing run (even to the same characters) forces the
# host to split that run at the sub-range boundaries, which is exactly
# the effect we want here.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Full text is "The picture first" (1-indexed characters):
#   1=T 2=h 3=e 4=' ' 5=p 6=i 7=c 8=t 9=u 10=r 11=e 12=' ' 13=f ...

# Split "The " -> "The" + " "
$tr.Characters(4, 1).Text = " "

# Split "picture " -> "picture" + " "
$tr.Characters(12, 1).Text = " "
